# Updates cryptocurrency Price (D) and Volume/1h (E) columns to the latest
# scraped figures. D-column cells hold text that sometimes *looks* numeric
# (e.g. "208.57"), so each of those is briefly switched to Text format before
# the value is written (and switched back to the default "Normal" style right
# after) so Excel keeps storing it as a string instead of auto-converting it to
# a number -- matching the original inlineStr cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cellRef -> new value
$updates = [ordered]@{
    'D2' = '27.064.86'
    'E2' = '  +0.56%  '
    'D3' = '1.566.07'
    'E3' = '  +1.05%  '
    'E4' = '  +0.56%  '
    'D5' = '208.57'
    'E5' = '  +1.07%  '
    'D6' = '0.491'
    'E6' = '  +0.81%  '
    'E7' = '  +0.47%  '
    'D8' = '22.10'
    'E8' = '  -0.01%  '
    'E9' = '  +0.98%  '
    'D10' = '0.0597'
    'E10' = '  +1.78%  '
    'D11' = '0.0859'
    'E11' = '  +0.45%  '
    'D12' = '1.568.51'
    'E12' = '  +1.52%  '
    'D13' = '3.79'
    'E13' = '  +1.40%  '
    'D14' = '0.520'
    'E14' = '  +0.35%  '
    'D15' = '27.049.49'
    'E15' = '  +0.52%  '
    'D16' = '61.87'
    'E16' = '  +0.39%  '
    'D17' = '0.0₃0706'
    'E17' = '  +1.17%  '
    'D18' = '7.43'
    'E18' = '  +2.28%  '
    'D19' = '215.79'
    'E19' = '  -0.51%  '
    'E20' = '  +0.45%  '
    'D21' = '4.15'
    'E21' = '  +2.41%  '
    'D22' = '9.20'
    'E22' = '  -0.27%  '
    'E23' = '  -0.12%  '
    'D24' = '154.08'
    'E24' = '  +0.08%  '
    'E25' = '  -0.27%  '
    'D26' = '15.06'
    'E26' = '  +0.74%  '
    'D27' = '0.105'
    'E27' = '  +1.44%  '
    'E28' = '  +0.45%  '
    'D29' = '0.0475'
    'E29' = '  +1.45%  '
    'E30' = '  +3.83%  '
    'E31' = '  +0.50%  '
    'D32' = '3.21'
    'E32' = '  +3.20%  '
    'D33' = '1.422.35'
    'E33' = '  +0.62%  '
    'D34' = '1.10'
    'E34' = '  +13.70%  '
    'E35' = '  +1.29%  '
    'D36' = '2.34'
    'E36' = '  +2.76%  '
    'D37' = '0.0166'
    'E37' = '  +0.90%  '
    'D38' = '0.534'
    'E38' = '  +1.48%  '
    'E39' = '  +2.36%  '
    'D40' = '0.811'
    'E40' = '  +0.39%  '
    'E41' = '  +0.44%  '
    'D42' = '2.34'
    'E42' = '  +0.61%  '
    'E43' = '  +0.47%  '
    'D44' = '64.74'
    'E44' = '  +0.20%  '
    'E45' = '  -0.42%  '
    'D46' = '1.703.98'
    'D47' = '86.79'
    'E47' = '  -0.71%  '
    'E48' = '  +3.52%  '
    'E49' = '  +0.62%  '
    'D50' = '0.0964'
    'E50' = '  +0.26%  '
    'E51' = '  +0.39%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $value = $updates[$ref]
    if ($ref.Substring(0,1) -eq 'D') {
        # Force text storage so numeric-looking prices stay strings.
        $cell.NumberFormat = '@'
        $cell.Value = $value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $value
    }
}
